$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 282.7879796666667
$ws.Cells.Item(2, 8).Value = 848.363939
$ws.Cells.Item(2, 9).Value = 0.9674521741401267
$ws.Cells.Item(2, 10).Value = 0.9674521741401266
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 7.214110666666667
$ws.Cells.Item(2, 14).Value = 21.642332
$ws.Cells.Item(2, 15).Value = 0.4688823795981188
$ws.Cells.Item(2, 16).Value = 0.4688823795981188
$ws.Cells.Item(2, 17).Value = 2040.063780518416
$ws.Cells.Item(2, 18).Value = 18360.57402466575
$ws.Cells.Item(2, 19).Value = 0.4536212775581963
$ws.Cells.Item(2, 20).Value = 0.4536212775581963

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 282.7879796666667
$ws.Cells.Item(3, 8).Value = 848.363939
$ws.Cells.Item(3, 9).Value = 0.9674521741401267
$ws.Cells.Item(3, 10).Value = 0.9674521741401266
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 7.110350666666666
$ws.Cells.Item(3, 14).Value = 21.331052
$ws.Cells.Item(3, 15).Value = 0.4621384803214003
$ws.Cells.Item(3, 16).Value = 0.4621384803214003
$ws.Cells.Item(3, 17).Value = 2010.721699748203
$ws.Cells.Item(3, 18).Value = 18096.49529773383
$ws.Cells.Item(3, 19).Value = 0.4470968775407529
$ws.Cells.Item(3, 20).Value = 0.4470968775407528

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 282.7879796666667
$ws.Cells.Item(4, 8).Value = 848.363939
$ws.Cells.Item(4, 9).Value = 0.9674521741401267
$ws.Cells.Item(4, 10).Value = 0.9674521741401266
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.061296333333333
$ws.Cells.Item(4, 14).Value = 3.183889
$ws.Cells.Item(4, 15).Value = 0.06897914008048092
$ws.Cells.Item(4, 16).Value = 0.06897914008048092
$ws.Cells.Item(4, 17).Value = 300.1218459309745
$ws.Cells.Item(4, 18).Value = 2701.096613378771
$ws.Cells.Item(4, 19).Value = 0.06673401904117762
$ws.Cells.Item(4, 20).Value = 0.06673401904117762

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 7.714696666666668
$ws.Cells.Item(5, 8).Value = 23.14409
$ws.Cells.Item(5, 9).Value = 0.02639291836872237
$ws.Cells.Item(5, 10).Value = 0.02639291836872237
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 7.214110666666667
$ws.Cells.Item(5, 14).Value = 21.642332
$ws.Cells.Item(5, 15).Value = 0.4688823795981188
$ws.Cells.Item(5, 16).Value = 0.4688823795981188
$ws.Cells.Item(5, 17).Value = 55.65467551309779
$ws.Cells.Item(5, 18).Value = 500.89207961788
$ws.Cells.Item(5, 19).Value = 0.01237517436926545
$ws.Cells.Item(5, 20).Value = 0.01237517436926545

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 7.714696666666668
$ws.Cells.Item(6, 8).Value = 23.14409
$ws.Cells.Item(6, 9).Value = 0.02639291836872237
$ws.Cells.Item(6, 10).Value = 0.02639291836872237
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 7.110350666666666
$ws.Cells.Item(6, 14).Value = 21.331052
$ws.Cells.Item(6, 15).Value = 0.4621384803214003
$ws.Cells.Item(6, 16).Value = 0.4621384803214003
$ws.Cells.Item(6, 17).Value = 54.85419858696445
$ws.Cells.Item(6, 18).Value = 493.6877872826801
$ws.Cells.Item(6, 19).Value = 0.01219718318616813
$ws.Cells.Item(6, 20).Value = 0.01219718318616813

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 7.714696666666668
$ws.Cells.Item(7, 8).Value = 23.14409
$ws.Cells.Item(7, 9).Value = 0.02639291836872237
$ws.Cells.Item(7, 10).Value = 0.02639291836872237
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.061296333333333
$ws.Cells.Item(7, 14).Value = 3.183889
$ws.Cells.Item(7, 15).Value = 0.06897914008048092
$ws.Cells.Item(7, 16).Value = 0.06897914008048092
$ws.Cells.Item(7, 17).Value = 8.187579285112221
$ws.Cells.Item(7, 18).Value = 73.68821356601
$ws.Cells.Item(7, 19).Value = 0.001820560813288799
$ws.Cells.Item(7, 20).Value = 0.001820560813288798

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.799090333333333
$ws.Cells.Item(8, 8).Value = 5.397271
$ws.Cells.Item(8, 9).Value = 0.006154907491150983
$ws.Cells.Item(8, 10).Value = 0.006154907491150983
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 7.214110666666667
$ws.Cells.Item(8, 14).Value = 21.642332
$ws.Cells.Item(8, 15).Value = 0.4688823795981188
$ws.Cells.Item(8, 16).Value = 0.4688823795981188
$ws.Cells.Item(8, 17).Value = 12.97883676399689
$ws.Cells.Item(8, 18).Value = 116.809530875972
$ws.Cells.Item(8, 19).Value = 0.00288592767065716
$ws.Cells.Item(8, 20).Value = 0.00288592767065716

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.799090333333333
$ws.Cells.Item(9, 8).Value = 5.397271
$ws.Cells.Item(9, 9).Value = 0.006154907491150983
$ws.Cells.Item(9, 10).Value = 0.006154907491150983
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 7.110350666666666
$ws.Cells.Item(9, 14).Value = 21.331052
$ws.Cells.Item(9, 15).Value = 0.4621384803214003
$ws.Cells.Item(9, 16).Value = 0.4621384803214003
$ws.Cells.Item(9, 17).Value = 12.79216315101022
$ws.Cells.Item(9, 18).Value = 115.129468359092
$ws.Cells.Item(9, 19).Value = 0.002844419594479318
$ws.Cells.Item(9, 20).Value = 0.002844419594479318

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.799090333333333
$ws.Cells.Item(10, 8).Value = 5.397271
$ws.Cells.Item(10, 9).Value = 0.006154907491150983
$ws.Cells.Item(10, 10).Value = 0.006154907491150983
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 1.061296333333333
$ws.Cells.Item(10, 14).Value = 3.183889
$ws.Cells.Item(10, 15).Value = 0.06897914008048092
$ws.Cells.Item(10, 16).Value = 0.06897914008048092
$ws.Cells.Item(10, 17).Value = 1.909367974102111
$ws.Cells.Item(10, 18).Value = 17.184311766919
$ws.Cells.Item(10, 19).Value = 0.0004245602260145051
$ws.Cells.Item(10, 20).Value = 0.0004245602260145051

